$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Through 2022-12-18 -> Through 2022-12-19)
$ws.Name = "Through 2022-12-19"

# Update header label in B1 (also updates shared string text)
$ws.Range("B1").Value = "December 2022 (through December 19)"

# Update counts across the sheet
$ws.Range("BV2").Value = 4
$ws.Range("Z5").Value = 2
$ws.Range("Z6").Value = 2
$ws.Range("Z7").Value = 5
$ws.Range("BJ9").Value = 5
$ws.Range("BJ14").Value = 6
$ws.Range("BV14").Value = 6
$ws.Range("CH14").Value = 3
$ws.Range("BJ20").Value = 5
$ws.Range("AL28").Value = 2
$ws.Range("B36").Value = 2
$ws.Range("Z39").Value = 2
$ws.Range("Z43").Value = 2
$ws.Range("B64").Value = 6
$ws.Range("N64").Value = 7
$ws.Range("BV65").Value = 2
$ws.Range("N72").Value = 1
$ws.Range("N83").Value = 3
$ws.Range("Z96").Value = 2
$ws.Range("B97").Value = 2
